$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Add a new "2022-Q1" sheet (positioned right after "2021-Q4", before the
#    "总计" summary sheet). We copy the "2021-Q4" sheet so the new sheet
#    inherits the identical header/column styling (bold header row + bordered,
#    bold column A), then we overwrite every cell with the 2022-Q1 data.
# ---------------------------------------------------------------------------
$srcSheet = $wb.Worksheets.Item("2021-Q4")
$srcSheet.Copy($null, $srcSheet)
$q1 = $wb.Worksheets.Item("2021-Q4 (2)")
$q1.Name = "2022-Q1"

# Extend column-A styling (bold + bordered, like A2/A3) down through row 12.
$q1.Range("A3").Copy()
$q1.Range("A4:A12").PasteSpecial(-4122)

# --- Header row ---
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# --- Data rows (A = index, B = 基金代码, C = 基金名称, D = 基金规模,
#     E = 股票总仓位, F = 仓位占比, G = 持有市值(亿元), H = 仓位排名) ---
# Numeric-looking text values (fund codes / percentages / amounts) are written
# with a leading apostrophe so Excel stores them as text, matching the source
# data (which keeps these as formatted strings rather than numbers).

$q1.Range("A2").Value = 0
$q1.Range("B2").Value = "'202027"
$q1.Range("C2").Value = "南方高端装备灵活配置混合A"
$q1.Range("D2").Value = "'15.67"
$q1.Range("E2").Value = "'90.90"
$q1.Range("F2").Value = "'2.90"
$q1.Range("G2").Value = "'0.4544"
$q1.Range("H2").Value = 9

$q1.Range("A3").Value = 1
$q1.Range("B3").Value = "'501030"
$q1.Range("C3").Value = "汇添富中证环境治理指数（LOF）A"
$q1.Range("D3").Value = "'6.61"
$q1.Range("E3").Value = "'93.20"
$q1.Range("F3").Value = "'2.66"
$q1.Range("G3").Value = "'0.1758"
$q1.Range("H3").Value = 1

$q1.Range("A4").Value = 2
$q1.Range("B4").Value = "'501031"
$q1.Range("C4").Value = "汇添富中证环境治理指数（LOF）C"
$q1.Range("D4").Value = "'2.74"
$q1.Range("E4").Value = "'93.20"
$q1.Range("F4").Value = "'2.66"
$q1.Range("G4").Value = "'0.0729"
$q1.Range("H4").Value = 1

$q1.Range("A5").Value = 3
$q1.Range("B5").Value = "'002295"
$q1.Range("C5").Value = "广发稳安灵活配置混合A"
$q1.Range("D5").Value = "'1.85"
$q1.Range("E5").Value = "'70.60"
$q1.Range("F5").Value = "'3.44"
$q1.Range("G5").Value = "'0.0636"
$q1.Range("H5").Value = 10

$q1.Range("A6").Value = 4
$q1.Range("B6").Value = "'005207"
$q1.Range("C6").Value = "南方高端装备灵活配置混合C"
$q1.Range("D6").Value = "'2.06"
$q1.Range("E6").Value = "'90.90"
$q1.Range("F6").Value = "'2.90"
$q1.Range("G6").Value = "'0.0597"
$q1.Range("H6").Value = 9

$q1.Range("A7").Value = 5
$q1.Range("B7").Value = "'005650"
$q1.Range("C7").Value = "万家量化同顺多策略灵活配置混合A"
$q1.Range("D7").Value = "'1.80"
$q1.Range("E7").Value = "'94.49"
$q1.Range("F7").Value = "'3.24"
$q1.Range("G7").Value = "'0.0583"
$q1.Range("H7").Value = 6

$q1.Range("A8").Value = 6
$q1.Range("B8").Value = "'164908"
$q1.Range("C8").Value = "交银施罗德中证环境治理指数（LOF）"
$q1.Range("D8").Value = "'2.12"
$q1.Range("E8").Value = "'93.72"
$q1.Range("F8").Value = "'2.41"
$q1.Range("G8").Value = "'0.0511"
$q1.Range("H8").Value = 2

$q1.Range("A9").Value = 7
$q1.Range("B9").Value = "'164401"
$q1.Range("C9").Value = "前海开源中证健康产业指数"
$q1.Range("D9").Value = "'2.13"
$q1.Range("E9").Value = "'94.15"
$q1.Range("F9").Value = "'1.25"
$q1.Range("G9").Value = "'0.0266"
$q1.Range("H9").Value = 6

$q1.Range("A10").Value = 8
$q1.Range("B10").Value = "'005651"
$q1.Range("C10").Value = "万家量化同顺多策略灵活配置混合C"
$q1.Range("D10").Value = "'0.29"
$q1.Range("E10").Value = "'94.49"
$q1.Range("F10").Value = "'3.24"
$q1.Range("G10").Value = "'0.0094"
$q1.Range("H10").Value = 6

$q1.Range("A11").Value = 9
$q1.Range("B11").Value = "'001657"
$q1.Range("C11").Value = "长安鑫富领先灵活配置混合"
$q1.Range("D11").Value = "'0.07"
$q1.Range("E11").Value = "'30.32"
$q1.Range("F11").Value = "'2.38"
$q1.Range("G11").Value = "'0.0017"
$q1.Range("H11").Value = 5

$q1.Range("A12").Value = 10
$q1.Range("B12").Value = "'008604"
$q1.Range("C12").Value = "广发稳安灵活配置混合C"
$q1.Range("D12").Value = "'0.02"
$q1.Range("E12").Value = "'70.60"
$q1.Range("F12").Value = "'3.44"
$q1.Range("G12").Value = "'0.0007"
$q1.Range("H12").Value = 10

# Drop the "text stored as number" quote-prefix styling that the leading
# apostrophe above implicitly applied, so these cells end up with the plain
# (unstyled) look of ordinary data cells - only column A keeps its bold /
# bordered style.
$q1.Range("B2:G12").Style = "Normal"

# ---------------------------------------------------------------------------
# 2) Update the "总计" (totals) sheet: insert the new 2022-Q1 summary row at
#    the top, pushing the existing quarters down by one row.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$b2 = $total.Range("B2").Value2
$c2 = $total.Range("C2").Value2
$d2 = $total.Range("D2").Value2

$b3 = $total.Range("B3").Value2
$c3 = $total.Range("C3").Value2
$d3 = $total.Range("D3").Value2

$b4 = $total.Range("B4").Value2
$c4 = $total.Range("C4").Value2
$d4 = $total.Range("D4").Value2

# Carry the bold/bordered column-A style down onto the new row 5.
$total.Range("A4").Copy()
$total.Range("A5").PasteSpecial(-4122)

$total.Range("A5").Value = 3
$total.Range("B5").Value = $b4
$total.Range("C5").Value = $c4
$total.Range("D5").Value = $d4

$total.Range("A4").Value = 2
$total.Range("B4").Value = $b3
$total.Range("C4").Value = $c3
$total.Range("D4").Value = $d3

$total.Range("A3").Value = 1
$total.Range("B3").Value = $b2
$total.Range("C3").Value = $c2
$total.Range("D3").Value = $d2

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 11
$total.Range("D2").Value = 0.97

# Restore the originally active sheet/selection (copying & renaming sheets
# along the way shifts focus onto the newest sheet).
$wb.Worksheets.Item("2021-Q1").Activate() | Out-Null
$wb.Worksheets.Item("2021-Q1").Range("A1").Select() | Out-Null

